$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9647635817527771
$ws.Range("B1").Value = 1.736628890037537
$ws.Range("C1").Value = 4.740855693817139
$ws.Range("D1").Value = 1.362407684326172
$ws.Range("E1").Value = 1.198545813560486
